$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.543.43', '  -0.67%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.852.13', '  -0.25%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9991', '  -0.15%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '243.38', '  -0.89%  ')
    ,@(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6366', '  -1.04%  ')
    ,@(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9999', '  -0.11%  ')
    ,@(8, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '48.30', '  +3.00%  ')
    ,@(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3005', '  +0.13%  ')
    ,@(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07479', '  -0.42%  ')
    ,@(11, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.29', '  +0.59%  ')
    ,@(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07628', '  -0.78%  ')
    ,@(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.908.66', '  +2.59%  ')
    ,@(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.033', '  -0.47%  ')
    ,@(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6867', '  +0.38%  ')
    ,@(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '83.65', '  -0.39%  ')
    ,@(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000009544', '  +0.90%  ')
    ,@(18, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.174', '  +1.62%  ')
    ,@(19, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.134.96', '  +0.62%  ')
    ,@(20, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.563.28', '  -0.59%  ')
    ,@(21, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '236.93', '  -1.52%  ')
    ,@(22, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.58', '  -0.96%  ')
    ,@(23, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.0000', '  -0.06%  ')
    ,@(24, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.740', '  +4.12%  ')
    ,@(25, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  -0.20%  ')
    ,@(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '157.35', '  -1.06%  ')
    ,@(27, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1406', '  -1.56%  ')
    ,@(28, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.509', '  -0.41%  ')
    ,@(29, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.77', '  -1.14%  ')
    ,@(30, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.492', '  -0.63%  ')
    ,@(31, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05988', '  -2.22%  ')
    ,@(32, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.250', '  -1.87%  ')
    ,@(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.121', '  -0.83%  ')
    ,@(34, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.073', '  -1.27%  ')
    ,@(35, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.873', '  -0.58%  ')
    ,@(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.176', '  +1.59%  ')
    ,@(37, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7204', '  -1.53%  ')
    ,@(38, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.605', '  -0.25%  ')
    ,@(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.792', '  -2.16%  ')
    ,@(40, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01776', '  -1.11%  ')
    ,@(41, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.202.67', '  -1.31%  ')
    ,@(42, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9095', '  -2.40%  ')
    ,@(43, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.161', '  -1.58%  ')
    ,@(44, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.052.41', '  +1.01%  ')
    ,@(45, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9995', '  -0.20%  ')
    ,@(46, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.92', '  -0.13%  ')
    ,@(47, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '66.80', '  +0.56%  ')
    ,@(48, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.297', '  +8.95%  ')
    ,@(49, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.00000000118', '  -4.23%  ')
    ,@(50, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4035', '  -1.20%  ')
    ,@(51, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.141', '  -1.53%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $bCell.NumberFormat = "@"
    $cCell.NumberFormat = "@"
    $dCell.NumberFormat = "@"
    $eCell.NumberFormat = "@"

    $bCell.Value = $row[1]
    $cCell.Value = $row[2]
    $dCell.Value = $row[3]
    $eCell.Value = $row[4]
}
